# Apply the edit described by the commit "All Result ready to start write":
#  1. Rename the shared string "MODEL_CONDITION" -> "MODELCONDITION"
#     (this is the header text in what is currently column E, row 1).
#  2. Remove the now-empty leading column A (the data previously lived in
#     columns B:F with column A only holding an unlabeled numeric/style
#     column) so everything shifts one column to the left (B:F -> A:E).
#     This also shrinks the used range from A1:F3 to A1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text first, while it is still easy to address by its
# current (pre-shift) column letter E.
$ws.Range("E1").Value = "MODELCONDITION"

# Delete column A entirely; this shifts B:F left into A:E and updates the
# sheet's used range/dimension accordingly.
$ws.Range("A:A").Delete()
